$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# The document currently ends with a (near) empty paragraph that only
# carries the _GoBack bookmark. The edit adds three brand-new
# paragraphs right before it, and appends one more run of text to
# that same bookmark paragraph (still in front of the bookmark).
#
# Word's Range.InsertXML, when given several <w:p> elements at once,
# creates standalone new paragraphs for every <w:p> except the last
# one, whose runs get merged into the paragraph that originally sat
# at the insertion point (preserving that paragraph's own pPr/rsids
# and, here, the bookmark). That is exactly the shape this diff wants,
# so we do it all in a single InsertXML call.
# ------------------------------------------------------------------

$n = $d.Paragraphs.Count
$tailPara = $d.Paragraphs.Item($n)
$insertPoint = $d.Range($tailPara.Range.Start, $tailPara.Range.Start)

$xmlAll = @"
<w:p $wns><w:r><w:t>Hvis tallet kan divideres med 5, er det enten 5, ellers er det sidste tal 5 eller 0. Derved kan man skrive:</w:t></w:r><w:r><w:br/><w:t>[5 | [0-9]*5 | [0-9]*0]</w:t></w:r></w:p><w:p $wns><w:r><w:t>Hvis tallet 5 skal opstå tre gange, kan tallet minimum være 555. Dog kan der også være et eller flere foran det første 5-tal, mellem de tre 5-taller og bag det sidste 5-tal.</w:t></w:r><w:r><w:br/><w:t>[555 | [0-46-9]*5</w:t></w:r><w:r><w:t>[0-46-9]*5[0-46-9]*5</w:t></w:r><w:r><w:t>[0-46-9]*]</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Overskrift2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>b)</w:t></w:r></w:p><w:p $wns><w:r><w:t>Det er et regulært sprog, hvis man kan bruge regular expression til at beskrive sproget. Dette kan godt lade sig gøre i det første eksampel (i), men kan ikke lade sig gøre i det andet eksempel (ii), eftersom man kun kan begrænse sig til en hvis længde i tallet, men tallet 1.000.000 er 7 lav, derfor vil man også kunne skrive alle tal op til 9.999.999, hvor der vil være tal som ikke indgår i det tiltænkte sprog.</w:t></w:r></w:p>
"@

$insertPoint.InsertXML($xmlAll) | Out-Null

Write-Output "done"
